$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.100864761350067
$ws.Range("C2").Value = 1.971249225817189
$ws.Range("D2").Value = 0.6998179082119407
$ws.Range("E2").Value = 0.2844193726699658
$ws.Range("G2").Value = 0.002691969173309583
$ws.Range("J2").Value = 0.147174348626514
$ws.Range("N2").Value = 3.949588442360181

$ws.Range("B3").Value = 5.866813803499497
$ws.Range("C3").Value = 1.89074214446714
$ws.Range("D3").Value = 0.6927438679833813
$ws.Range("E3").Value = 0.2807244758264034
$ws.Range("G3").Value = 0.002702811064109104
$ws.Range("J3").Value = 0.1445661577799342
$ws.Range("N3").Value = 3.866634005347549

$ws.Range("B4").Value = 5.727770283480652
$ws.Range("C4").Value = 1.842860936349553
$ws.Range("D4").Value = 0.6888396703855904
$ws.Range("E4").Value = 0.2786431498179809
$ws.Range("G4").Value = 0.002709792206916302
$ws.Range("J4").Value = 0.1430681449162847
$ws.Range("N4").Value = 3.816141788802298

$ws.Range("B5").Value = 5.672261971899104
$ws.Range("C5").Value = 1.823731509942206
$ws.Range("D5").Value = 0.6873581605358936
$ws.Range("E5").Value = 0.2778416245307653
$ws.Range("G5").Value = 0.002712719017721501
$ws.Range("J5").Value = 0.1424833827619878
$ws.Range("N5").Value = 3.795671840631087

$ws.Range("B6").Value = 5.663113927774248
$ws.Range("C6").Value = 1.820577978083577
$ws.Range("D6").Value = 0.6871187413045448
$ws.Range("E6").Value = 0.2777113345509719
$ws.Range("G6").Value = 0.002713209973696837
$ws.Range("J6").Value = 0.1423878259690383
$ws.Range("N6").Value = 3.792279077129479

$ws.Range("B7").Value = 5.72701703625944
$ws.Range("C7").Value = 1.842601410802445
$ws.Range("D7").Value = 0.6888192482079489
$ws.Range("E7").Value = 0.2786321519799202
$ws.Range("G7").Value = 0.002709831346625884
$ws.Range("J7").Value = 0.143060155002587
$ws.Range("N7").Value = 3.815865301308406

$ws.Range("B8").Value = 6.019183598494237
$ws.Range("C8").Value = 1.9431640945769
$ws.Range("D8").Value = 0.6972870009137182
$ws.Range("E8").Value = 0.2831061916559747
$ws.Range("G8").Value = 0.002695640428933691
$ws.Range("J8").Value = 0.1462533834890962
$ws.Range("N8").Value = 3.920890844753018

$ws.Range("B9").Value = 6.630089918603744
$ws.Range("C9").Value = 2.153019500780772
$ws.Range("D9").Value = 0.7174250967493663
$ws.Range("E9").Value = 0.2933894122547542
$ws.Range("G9").Value = 0.002670364924978603
$ws.Range("J9").Value = 0.1533511226293029
$ws.Range("N9").Value = 4.130598401389221

$ws.Range("B10").Value = 7.103478902717484
$ws.Range("C10").Value = 2.315427690123556
$ws.Range("D10").Value = 0.7344442743708726
$ws.Range("E10").Value = 0.3018988481844644
$ws.Range("G10").Value = 0.002653323938024032
$ws.Range("J10").Value = 0.1590977043080812
$ws.Range("N10").Value = 4.287315866659924

$ws.Range("B11").Value = 7.324488363274099
$ws.Range("C11").Value = 2.391214666883229
$ws.Range("D11").Value = 0.7426856237164827
$ws.Range("E11").Value = 0.305985014213114
$ws.Range("G11").Value = 0.002645897585273005
$ws.Range("J11").Value = 0.1618325515300114
$ws.Range("N11").Value = 4.359265155112496

$ws.Range("B12").Value = 7.409019366831672
$ws.Range("C12").Value = 2.420196985552366
$ws.Range("D12").Value = 0.7458795231767681
$ws.Range("E12").Value = 0.3075639241097434
$ws.Range("G12").Value = 0.002643131777951162
$ws.Range("J12").Value = 0.1628859434558194
$ws.Range("N12").Value = 4.386611467176579

$ws.Range("B13").Value = 7.390776353333536
$ws.Range("C13").Value = 2.413942367992718
$ws.Range("D13").Value = 0.7451883885012194
$ws.Range("E13").Value = 0.3072224634425993
$ws.Range("G13").Value = 0.002643725387635756
$ws.Range("J13").Value = 0.1626582799103176
$ws.Range("N13").Value = 4.380717373164316

$ws.Range("B14").Value = 7.331425806803736
$ws.Range("C14").Value = 2.393593325043753
$ws.Range("D14").Value = 0.7429469155655966
$ws.Range("E14").Value = 0.3061142757577286
$ws.Range("G14").Value = 0.002645669113170815
$ws.Range("J14").Value = 0.1619188562950598
$ws.Range("N14").Value = 4.36151289946946

$ws.Range("B15").Value = 7.29518197804623
$ws.Range("C15").Value = 2.381166140513301
$ws.Range("D15").Value = 0.7415835048661847
$ws.Range("E15").Value = 0.3054396087647007
$ws.Range("G15").Value = 0.002646865729486424
$ws.Range("J15").Value = 0.1614682638334557
$ws.Range("N15").Value = 4.349762899625489

$ws.Range("B16").Value = 7.089151497454054
$ws.Range("C16").Value = 2.310513984355453
$ws.Range("D16").Value = 0.7339158378822503
$ws.Range("E16").Value = 0.3016361904935962
$ws.Range("G16").Value = 0.00265381578339587
$ws.Range("J16").Value = 0.1589214389446028
$ws.Range("N16").Value = 4.282627564391589

$ws.Range("B17").Value = 6.964226461472663
$ws.Range("C17").Value = 2.26766600856331
$ws.Range("D17").Value = 0.7293407417625417
$ws.Range("E17").Value = 0.2993584668247991
$ws.Range("G17").Value = 0.002658162518469226
$ws.Range("J17").Value = 0.1573902440723316
$ws.Range("N17").Value = 4.241614961437364

$ws.Range("B18").Value = 6.892904297272707
$ws.Range("C18").Value = 2.243199879084273
$ws.Range("D18").Value = 0.7267561782745702
$ws.Range("E18").Value = 0.2980685955132856
$ws.Range("G18").Value = 0.00266069332635093
$ws.Range("J18").Value = 0.1565208766863293
$ws.Range("N18").Value = 4.218087131693181

$ws.Range("B19").Value = 6.868846352777382
$ws.Range("C19").Value = 2.234946516096556
$ws.Range("D19").Value = 0.7258891079746945
$ws.Range("E19").Value = 0.2976353205523665
$ws.Range("G19").Value = 0.002661555497407255
$ws.Range("J19").Value = 0.156228458488215
$ws.Range("N19").Value = 4.21013139812834

$ws.Range("B20").Value = 6.977469765405772
$ws.Range("C20").Value = 2.272208667943346
$ws.Range("D20").Value = 0.7298229042852711
$ws.Range("E20").Value = 0.2995988377604348
$ws.Range("G20").Value = 0.002657696629005762
$ws.Range("J20").Value = 0.1575520662979102
$ws.Range("N20").Value = 4.245974410738853

$ws.Range("B21").Value = 7.348835511029506
$ws.Range("C21").Value = 2.399562564524444
$ws.Range("D21").Value = 0.743603296629999
$ws.Range("E21").Value = 0.3064389155168428
$ws.Range("G21").Value = 0.002645096938412635
$ws.Range("J21").Value = 0.1621355574804824
$ws.Range("N21").Value = 4.367150937174586

$ws.Range("B22").Value = 7.596450200032677
$ws.Range("C22").Value = 2.484452301861097
$ws.Range("D22").Value = 0.7530361759088464
$ws.Range("E22").Value = 0.311093608081876
$ws.Range("G22").Value = 0.002637132509443156
$ws.Range("J22").Value = 0.1652348799987493
$ws.Range("N22").Value = 4.446936571616618

$ws.Range("B23").Value = 7.463836130534105
$ws.Range("C23").Value = 2.438990329929652
$ws.Range("D23").Value = 0.7479621981686364
$ws.Range("E23").Value = 0.308592234474574
$ws.Range("G23").Value = 0.002641358698740155
$ws.Range("J23").Value = 0.1635710816721314
$ws.Range("N23").Value = 4.404297457991902

$ws.Range("B24").Value = 6.971480915626216
$ws.Range("C24").Value = 2.270154409854058
$ws.Range("D24").Value = 0.7296047762975491
$ws.Range("E24").Value = 0.2994901049880667
$ws.Range("G24").Value = 0.002657907158504068
$ws.Range("J24").Value = 0.1574788724343676
$ws.Range("N24").Value = 4.244003345745085

$ws.Range("B25").Value = 6.460615802797065
$ws.Range("C25").Value = 2.094841366924925
$ws.Range("D25").Value = 0.7115919104017792
$ws.Range("E25").Value = 0.2904424652754827
$ws.Range("G25").Value = 0.002676932132875486
$ws.Range("J25").Value = 0.1513392525463289
$ws.Range("N25").Value = 4.073428694324036
